# Update "In Class Demonstrations" sheet:
#  - Update regression coefficients in row 2 (U Lag) and row 3 (C/A Lag)
#  - Remove the "Constant" row (row 4) and "r2_adj" row (row 5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format so numeric-looking values stay as text (shared strings),
# matching the source data which stores these as text, not numbers.
$ws.Range("B2:C3").NumberFormat = "@"

# Update the changed coefficient values (keep text formatting, e.g. "-0.01" not "-0.009")
# Assign column by column (B2, B3, then C2, C3) to match original authoring order.
$ws.Range("B2").Value = "-0.358***"
$ws.Range("B3").Value = "-0.01"
$ws.Range("C2").Value = "1.248"
$ws.Range("C3").Value = "-0.351***"

# Restore the default (General) styling so the cells match the original
# workbook's default style (no explicit style index).
$ws.Range("B2:C3").ClearFormats()

# Remove the last two rows (Constant and r2_adj) entirely
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
